$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1234.1666
$ws.Range("J17").Value = 1025.0426
$ws.Range("L17").Value = 3075.1278
$ws.Range("N17").Value = -3411.1278

$ws.Range("H58").Value = 1155.3846
$ws.Range("I58").Value = 324.44446
$ws.Range("J58").Value = 3025
$ws.Range("K58").Value = 973.33338
$ws.Range("L58").Value = 9075
$ws.Range("M58").Value = -823.33338
$ws.Range("N58").Value = -9375

$ws.Range("H62").Value = 6139.143
$ws.Range("I62").Value = 999.6667
$ws.Range("K62").Value = 999.6667
$ws.Range("M62").Value = -375.6667

$ws.Range("H65").Value = 6139.143
$ws.Range("I65").Value = 999.6667
$ws.Range("K65").Value = 4998.3335
$ws.Range("M65").Value = -1878.3335

$ws.Range("H129").Value = 888.7183
$ws.Range("I129").Value = 999.8570999999999
$ws.Range("J129").Value = 876.5625
$ws.Range("K129").Value = 2999.5713
$ws.Range("L129").Value = 2629.6875
$ws.Range("M129").Value = 2000.4287
$ws.Range("N129").Value = -12629.6875

$ws.Range("H131").Value = 2892.0715
$ws.Range("J131").Value = 4427.5
$ws.Range("L131").Value = 13282.5
$ws.Range("N131").Value = -23362.5

$ws.Range("H132").Value = 1105.1724
$ws.Range("I132").Value = 1117.6923
$ws.Range("K132").Value = 3353.0769
$ws.Range("M132").Value = -823.0769

$ws.Range("H137").Value = 2103.6667
$ws.Range("I137").Value = 1238.2
$ws.Range("J137").Value = 2721.8572
$ws.Range("K137").Value = 3714.6
$ws.Range("L137").Value = 8165.571599999999
$ws.Range("M137").Value = -1164.6
$ws.Range("N137").Value = -13265.5716

$ws.Range("H139").Value = 46999.5
$ws.Range("J139").Value = 46999.5
$ws.Range("L139").Value = 46999.5
$ws.Range("N139").Value = -57279.5

$ws.Range("H141").Value = 3978.4666
$ws.Range("I141").Value = 3349.6667
$ws.Range("J141").Value = 4397.6665
$ws.Range("K141").Value = 10049.0001
$ws.Range("L141").Value = 13192.9995
$ws.Range("M141").Value = -4869.000100000001
$ws.Range("N141").Value = -23552.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3155.4558
$ws.Range("I32").Value = 1641.4822
$ws.Range("K32").Value = 1641.4822
$ws.Range("M32").Value = -1354.4822

$ws.Range("H45").Value = 1423.8948
$ws.Range("I45").Value = 1008
$ws.Range("K45").Value = 1008
$ws.Range("M45").Value = -631

$ws.Range("H82").Value = 83333
$ws.Range("J82").Value = 83333
$ws.Range("L82").Value = 83333
$ws.Range("N82").Value = -84055

$ws.Range("H85").Value = 83333
$ws.Range("J85").Value = 83333
$ws.Range("L85").Value = 83333
$ws.Range("N85").Value = -85829

$ws.Range("H132").Value = 1648.9412
$ws.Range("I132").Value = 1288.0714
$ws.Range("K132").Value = 3864.2142
$ws.Range("M132").Value = -1334.2142

$ws.Range("H135").Value = 44614.5
$ws.Range("J135").Value = 44614.5
$ws.Range("L135").Value = 44614.5
$ws.Range("N135").Value = -54754.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 6452
$ws.Range("I80").Value = 41.5
$ws.Range("K80").Value = 41.5
$ws.Range("M80").Value = 956.5

$ws.Range("H83").Value = 6452
$ws.Range("I83").Value = 41.5
$ws.Range("K83").Value = 207.5
$ws.Range("M83").Value = 4784.5

$ws.Range("H105").Value = 1712.8518
$ws.Range("I105").Value = 1740.3077
$ws.Range("K105").Value = 1740.3077
$ws.Range("M105").Value = 6.692299999999932

$ws.Range("H134").Value = 20480
$ws.Range("I134").Value = 27865.143
$ws.Range("J134").Value = 10140.8
$ws.Range("K134").Value = 83595.429
$ws.Range("L134").Value = 30422.4
$ws.Range("M134").Value = -81060.429
$ws.Range("N134").Value = -35492.39999999999

$ws.Range("H135").Value = 55390
$ws.Range("J135").Value = 55390
$ws.Range("L135").Value = 55390
$ws.Range("N135").Value = -65530

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 659.8
$ws.Range("I16").Value = 576.5
$ws.Range("K16").Value = 576.5
$ws.Range("M16").Value = -289.5

$ws.Range("H31").Value = 2966.2693
$ws.Range("I31").Value = 1965.45
$ws.Range("J31").Value = 6302.3335
$ws.Range("K31").Value = 1965.45
$ws.Range("L31").Value = 6302.3335
$ws.Range("M31").Value = -1670.45
$ws.Range("N31").Value = -6892.3335

$ws.Range("H34").Value = 2966.2693
$ws.Range("I34").Value = 1965.45
$ws.Range("J34").Value = 6302.3335
$ws.Range("K34").Value = 1965.45
$ws.Range("L34").Value = 6302.3335
$ws.Range("M34").Value = -1763.45
$ws.Range("N34").Value = -6706.3335

$ws.Range("H113").Value = 659.8
$ws.Range("I113").Value = 576.5
$ws.Range("K113").Value = 576.5
$ws.Range("M113").Value = 1593.5

$ws.Range("H134").Value = 1865.8889
$ws.Range("I134").Value = 1436.6875
$ws.Range("J134").Value = 5299.5
$ws.Range("K134").Value = 4310.0625
$ws.Range("L134").Value = 15898.5
$ws.Range("M134").Value = -1775.0625
$ws.Range("N134").Value = -20968.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 693.4286
$ws.Range("I5").Value = 670.8
$ws.Range("K5").Value = 2012.4
$ws.Range("M5").Value = -1900.4

$ws.Range("H51").Value = 900
$ws.Range("I51").Value = 900
$ws.Range("K51").Value = 2700
$ws.Range("M51").Value = -2240

$ws.Range("H68").Value = 800
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 800
$ws.Range("K68").Value = 0
$ws.Range("N68").Value = -4022
$ws.Range("M68").ClearContents()

$ws.Range("H71").Value = 800
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 800
$ws.Range("K71").Value = 0
$ws.Range("N71").Value = -15312
$ws.Range("M71").ClearContents()

$ws.Range("H103").Value = 2158.6
$ws.Range("I103").Value = 1720.6
$ws.Range("J103").Value = 2596.6
$ws.Range("K103").Value = 5161.799999999999
$ws.Range("L103").Value = 7789.799999999999
$ws.Range("M103").Value = -4282.799999999999
$ws.Range("N103").Value = -9547.799999999999

$ws.Range("H113").Value = 8490.846
$ws.Range("I113").Value = 50352
$ws.Range("J113").Value = 879.7273
$ws.Range("K113").Value = 151056
$ws.Range("L113").Value = 2639.1819
$ws.Range("M113").Value = -148886
$ws.Range("N113").Value = -6979.1819

$ws.Range("H131").Value = 769.7071
$ws.Range("J131").Value = 788.70966
$ws.Range("L131").Value = 2366.12898
$ws.Range("N131").Value = -12446.12898

$ws.Range("H132").Value = 1792.5
$ws.Range("J132").Value = 1811
$ws.Range("L132").Value = 16299
$ws.Range("N132").Value = -21359

$ws.Range("H134").Value = 1801.0416
$ws.Range("I134").Value = 1361.65
$ws.Range("K134").Value = 4084.95
$ws.Range("M134").Value = 985.0499999999997

$ws.Range("H135").Value = 693.4286
$ws.Range("I135").Value = 670.8
$ws.Range("K135").Value = 6037.2
$ws.Range("M135").Value = -3502.2

$ws.Range("H140").Value = 1742.72
$ws.Range("I140").Value = 864.8
$ws.Range("J140").Value = 2328
$ws.Range("K140").Value = 2594.4
$ws.Range("L140").Value = 6984
$ws.Range("M140").Value = 2585.6
$ws.Range("N140").Value = -17344

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2029037.4
$ws.Range("J132").Value = 6593.3335
$ws.Range("L132").Value = 19780.0005
$ws.Range("N132").Value = -24840.0005

$ws.Range("H140").Value = 48333
$ws.Range("J140").Value = 48333
$ws.Range("L140").Value = 48333
$ws.Range("N140").Value = -58693

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3081.2856
$ws.Range("J22").Value = 2113.8
$ws.Range("L22").Value = 2113.8
$ws.Range("N22").Value = -2703.8

$ws.Range("H27").Value = 3081.2856
$ws.Range("J27").Value = 2113.8
$ws.Range("L27").Value = 2113.8
$ws.Range("N27").Value = -2327.8

$ws.Range("H61").Value = 3474.625
$ws.Range("I61").Value = 3200
$ws.Range("J61").Value = 3749.25
$ws.Range("K61").Value = 3200
$ws.Range("L61").Value = 3749.25
$ws.Range("M61").Value = -2998
$ws.Range("N61").Value = -4153.25

$ws.Range("H113").Value = 3474.625
$ws.Range("I113").Value = 3200
$ws.Range("J113").Value = 3749.25
$ws.Range("K113").Value = 3200
$ws.Range("L113").Value = 3749.25
$ws.Range("M113").Value = -1030
$ws.Range("N113").Value = -8089.25

$ws.Range("H132").Value = 2322.5
$ws.Range("I132").Value = 2532.3333
$ws.Range("J132").Value = 2265.2727
$ws.Range("K132").Value = 7596.999899999999
$ws.Range("L132").Value = 6795.8181
$ws.Range("M132").Value = -5066.999899999999
$ws.Range("N132").Value = -11855.8181

$ws.Range("H133").Value = 79163
$ws.Range("J133").Value = 79163
$ws.Range("L133").Value = 79163
$ws.Range("N133").Value = -84223

$ws.Range("H134").Value = 48808.8
$ws.Range("J134").Value = 48808.8
$ws.Range("L134").Value = 48808.8
$ws.Range("N134").Value = -58948.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 6497.722
$ws.Range("I132").Value = 1928.2
$ws.Range("K132").Value = 5784.6
$ws.Range("M132").Value = -3254.6
